# Update the task_point data values on the active sheet.
# The diff only changes the numeric values in A1:F10; everything else
# (layout, formatting, dimensions) stays the same.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(-0.097, -0.168, -0.138, 36, -46, -44),
    @(-0.133, 0.08699999999999999, -0.081, 96, -32, 158),
    @(0.031, -0.281, -0.159, -90, 28, -62),
    @(-0.029, 0.174, 0.002, -91, 35, -3),
    @(0.068, -0.115, -0.193, 14, 61, 35),
    @(0.104, 0.281, -0.392, -62, -18, -31),
    @(-0.261, 0.01, -0.108, 86, 45, 4),
    @(0.077, -0.14, 0.152, -64, 3, 175),
    @(0.014, 0.336, -0.041, -89, 18, 94),
    @(-0.143, -0.168, -0.141, -166, -19, 73)
)

for ($r = 1; $r -le 10; $r++) {
    for ($c = 1; $c -le 6; $c++) {
        $ws.Cells.Item($r, $c).Value = $data[$r - 1][$c - 1]
    }
}
